# Breadcrumbs GSC export refresh ("updated all GSC export files"): a new
# day of data ("2025-11-07") became available, so it is appended as a new
# row at the bottom of the "Chart" sheet, exactly like every prior daily
# row (date label in col A, the Invalid/Valid counters in cols B/C).
# "Critical issues" / "Non-critical issues" already hold the right
# Issue/Validation/Items header text and need no further changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.UsedRange.Row() + $ws.UsedRange.Rows.Count() - 1
$newRow = $lastRow + 1

# Clone the cell directly above first: this carries over its "plain text"
# shared-string cell type (rather than letting the Value setter
# auto-detect a date-like string and convert it into a date
# serial-number/date-formatted cell -- which is NOT how the existing date
# column is stored; every prior row keeps its date as plain text).
$dateCell = $ws.Range("A" + $newRow)
$ws.Range("A" + $lastRow).Copy($dateCell)

# Swap in the new date text without re-triggering that autodetection:
# build it as a literal-string formula, then paste-special the computed
# value back over itself so the final cell is plain text again (no
# formula left behind, no number-format override needed).
$dateCell.Formula = '="2025-11-07"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B" + $newRow).Value = 0
$ws.Range("C" + $newRow).Value = 87
